$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple text / non-numeric-looking values (safe to assign directly)
$ws.Range("D2").Value = "33.714.49"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.765.31"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +3.05%  "
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").Value = "2.018.53"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("E13").Value = "  +7.35%  "
$ws.Range("D14").Value = "1.751.13"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "33.699.58"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("E32").Value = "  -3.13%  "
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").Value = "1.380.06"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  +5.47%  "
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E43").Value = "  +15.68%  "
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("E45").Value = "  +15.75%  "
$ws.Range("E46").Value = "  +4.80%  "
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").Value = "1.920.11"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  +0.43%  "

# Numeric-looking strings that must remain TEXT (force text format, set, then clear format)
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "224.25"
$c.ClearFormats()
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "31.96"
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "4.14"
$c.ClearFormats()
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "66.53"
$c.ClearFormats()
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "237.50"
$c.ClearFormats()
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.59"
$c.ClearFormats()
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "159.45"
$c.ClearFormats()
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "16.11"
$c.ClearFormats()
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.03"
$c.ClearFormats()
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0510"
$c.ClearFormats()
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.36"
$c.ClearFormats()
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.907"
$c.ClearFormats()
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "77.63"
$c.ClearFormats()
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "13.53"
$c.ClearFormats()
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "107.69"
$c.ClearFormats()
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "5.81"
$c.ClearFormats()
